$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2: Target cluster changes from MuSCs -> ECs,
# and the per-edge stats are recomputed for the new ("MuSCs"->"ECs") pairing.
$ws.Range("D2").Value = "ECs"

$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.004083333333333334
$ws.Range("N2").Value = 0.01225
$ws.Range("O2").Value = 0.1210581968752162
$ws.Range("P2").Value = 0.1210581968752162
$ws.Range("Q2").Value = 0.00001285161111111111
$ws.Range("R2").Value = 0.0001156645
$ws.Range("S2").Value = 0.1210581968752162
$ws.Range("T2").Value = 0.1210581968752162

# --- Insert a new row 3 holding the original MuSCs -> MuSCs pairing
# (the stats that used to live on row 2 before the TPM recompute), with
# the last four ratio columns (O,P,S,T) updated for the new total.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Dsc3"
$ws.Range("C3").Value = "Dsg1a"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.003147333333333334
$ws.Range("H3").Value = 0.009442000000000001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.029647
$ws.Range("N3").Value = 0.08894100000000001
$ws.Range("O3").Value = 0.8789418031247839
$ws.Range("P3").Value = 0.8789418031247839
$ws.Range("Q3").Value = 0.00009330899133333336
$ws.Range("R3").Value = 0.0008397809220000002
$ws.Range("S3").Value = 0.8789418031247839
$ws.Range("T3").Value = 0.8789418031247839
